$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 6399.4165
$ws.Range("I32").Value = 10799.5
$ws.Range("K32").Value = 10799.5
$ws.Range("M32").Value = -10473.5
$ws.Range("H40").Value = 3175
$ws.Range("J40").Value = 1900
$ws.Range("L40").Value = 1900
$ws.Range("N40").Value = -2250
$ws.Range("H51").Value = 15000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 45459130
$ws.Range("I132").Value = 50004816
$ws.Range("J132").Value = 2246
$ws.Range("K132").Value = 150014448
$ws.Range("L132").Value = 6738
$ws.Range("M132").Value = -150011918
$ws.Range("N132").Value = -11798
$ws.Range("H141").Value = 2654.4546
$ws.Range("I141").Value = 2274.625
$ws.Range("K141").Value = 6823.875
$ws.Range("M141").Value = -1643.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2842.2727
$ws.Range("I2").Value = 2461.3333
$ws.Range("K2").Value = 2461.3333
$ws.Range("M2").Value = -2348.3333
$ws.Range("H32").Value = 8460.312
$ws.Range("I32").Value = 7741.122
$ws.Range("K32").Value = 7741.122
$ws.Range("M32").Value = -7454.122
$ws.Range("H43").Value = 43333.332
$ws.Range("H45").Value = 1538.8
$ws.Range("I45").Value = 1298.5
$ws.Range("K45").Value = 1298.5
$ws.Range("M45").Value = -921.5
$ws.Range("H56").Value = 7000
$ws.Range("I56").Value = 7000
$ws.Range("K56").Value = 7000
$ws.Range("M56").Value = -6258
$ws.Range("H74").Value = 2168.182
$ws.Range("I74").Value = 1886.1
$ws.Range("J74").Value = 4989
$ws.Range("K74").Value = 1886.1
$ws.Range("L74").Value = 4989
$ws.Range("M74").Value = -1012.1
$ws.Range("N74").Value = -6737
$ws.Range("H77").Value = 2168.182
$ws.Range("I77").Value = 1886.1
$ws.Range("J77").Value = 4989
$ws.Range("K77").Value = 9430.5
$ws.Range("L77").Value = 24945
$ws.Range("M77").Value = -5062.5
$ws.Range("N77").Value = -33681
$ws.Range("H116").Value = 2842.2727
$ws.Range("I116").Value = 2461.3333
$ws.Range("K116").Value = 2461.3333
$ws.Range("M116").Value = -167.3332999999998
$ws.Range("H132").Value = 4509.2856
$ws.Range("I132").Value = 4509.2856
$ws.Range("K132").Value = 13527.8568
$ws.Range("M132").Value = -10997.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2842.2727
$ws.Range("I3").Value = 2461.3333
$ws.Range("K3").Value = 2461.3333
$ws.Range("M3").Value = -2347.3333
$ws.Range("H86").Value = 18560508
$ws.Range("I86").Value = 47935.273
$ws.Range("J86").Value = 47651692
$ws.Range("K86").Value = 47935.273
$ws.Range("L86").Value = 47651692
$ws.Range("M86").Value = -46812.273
$ws.Range("N86").Value = -47653938
$ws.Range("H89").Value = 18560508
$ws.Range("I89").Value = 47935.273
$ws.Range("J89").Value = 47651692
$ws.Range("K89").Value = 239676.365
$ws.Range("L89").Value = 238258460
$ws.Range("M89").Value = -234060.365
$ws.Range("N89").Value = -238269692

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8320.091
$ws.Range("I31").Value = 8591.223
$ws.Range("J31").Value = 7100
$ws.Range("K31").Value = 8591.223
$ws.Range("L31").Value = 7100
$ws.Range("M31").Value = -8296.223
$ws.Range("N31").Value = -7690
$ws.Range("H34").Value = 8320.091
$ws.Range("I34").Value = 8591.223
$ws.Range("J34").Value = 7100
$ws.Range("K34").Value = 8591.223
$ws.Range("L34").Value = 7100
$ws.Range("M34").Value = -8389.223
$ws.Range("N34").Value = -7504
$ws.Range("H134").Value = 1976.6333
$ws.Range("I134").Value = 2067.4167
$ws.Range("J134").Value = 1613.5
$ws.Range("K134").Value = 6202.250100000001
$ws.Range("L134").Value = 4840.5
$ws.Range("M134").Value = -3667.250100000001
$ws.Range("N134").Value = -9910.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2524058
$ws.Range("I4").Value = 28849.6
$ws.Range("J4").Value = 15000100
$ws.Range("K4").Value = 86548.79999999999
$ws.Range("L4").Value = 45000300
$ws.Range("M4").Value = -86436.79999999999
$ws.Range("N4").Value = -45000524
$ws.Range("H19").Value = 4900
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 4900
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 14700
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -15048
$ws.Range("H34").Value = 49543.305
$ws.Range("J34").Value = 56921.25
$ws.Range("L34").Value = 170763.75
$ws.Range("N34").Value = -170931.75
$ws.Range("H62").Value = 9845.200000000001
$ws.Range("I62").Value = 2114.5
$ws.Range("J62").Value = 14999
$ws.Range("K62").Value = 6343.5
$ws.Range("L62").Value = 44997
$ws.Range("M62").Value = -5657.5
$ws.Range("N62").Value = -46369
$ws.Range("H65").Value = 9845.200000000001
$ws.Range("I65").Value = 2114.5
$ws.Range("J65").Value = 14999
$ws.Range("K65").Value = 19030.5
$ws.Range("L65").Value = 134991
$ws.Range("M65").Value = -15598.5
$ws.Range("N65").Value = -141855
$ws.Range("H131").Value = 1977.25
$ws.Range("I131").Value = 1478.4286
$ws.Range("J131").Value = 2675.6
$ws.Range("K131").Value = 4435.2858
$ws.Range("L131").Value = 8026.799999999999
$ws.Range("M131").Value = 604.7142000000003
$ws.Range("N131").Value = -18106.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1385.7273
$ws.Range("I102").Value = 1392.4138
$ws.Range("K102").Value = 1392.4138
$ws.Range("M102").Value = 229.5862
$ws.Range("H122").Value = 2251.25
$ws.Range("I122").Value = 2335
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 7005
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -4555
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H62").Value = 28666.666
$ws.Range("J62").Value = 27500
$ws.Range("L62").Value = 27500
$ws.Range("N62").Value = -28748
$ws.Range("H64").Value = 49999.5
$ws.Range("J64").Value = 49999.5
$ws.Range("L64").Value = 49999.5
$ws.Range("N64").Value = -50449.5
$ws.Range("H65").Value = 28666.666
$ws.Range("J65").Value = 27500
$ws.Range("L65").Value = 82500
$ws.Range("N65").Value = -88740
$ws.Range("H67").Value = 49999.5
$ws.Range("J67").Value = 49999.5
$ws.Range("L67").Value = 49999.5
$ws.Range("N67").Value = -51559.5
$ws.Range("H122").Value = 2964.1428
$ws.Range("I122").Value = 2670
$ws.Range("J122").Value = 3699.5
$ws.Range("K122").Value = 8010
$ws.Range("L122").Value = 11098.5
$ws.Range("M122").Value = -5560
$ws.Range("N122").Value = -15998.5
$ws.Range("H132").Value = 3416.3333
$ws.Range("I132").Value = 3416.3333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10248.9999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7718.999899999999
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H30").Value = 16993.25
$ws.Range("I30").Value = 19325
$ws.Range("J30").Value = 9998
$ws.Range("K30").Value = 19325
$ws.Range("L30").Value = 9998
$ws.Range("M30").Value = -19218
$ws.Range("N30").Value = -10212
$ws.Range("H46").Value = 47500
$ws.Range("J46").Value = 47500
$ws.Range("L46").Value = 47500
$ws.Range("N46").Value = -47962
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H63").Value = 48914.668
$ws.Range("J63").Value = 48914.668
$ws.Range("L63").Value = 48914.668
$ws.Range("N63").Value = -50162.668
$ws.Range("H66").Value = 48914.668
$ws.Range("J66").Value = 48914.668
$ws.Range("L66").Value = 146744.004
$ws.Range("N66").Value = -152984.004
$ws.Range("H81").Value = 8558.223
$ws.Range("I81").Value = 4822.1113
$ws.Range("K81").Value = 9644.222599999999
$ws.Range("M81").Value = -8583.222599999999
$ws.Range("H82").Value = 70272.5
$ws.Range("J82").Value = 70272
$ws.Range("L82").Value = 70272
$ws.Range("N82").Value = -71038
$ws.Range("H84").Value = 8558.223
$ws.Range("I84").Value = 4822.1113
$ws.Range("K84").Value = 48221.113
$ws.Range("M84").Value = -42917.113
$ws.Range("H85").Value = 70272.5
$ws.Range("J85").Value = 70272
$ws.Range("L85").Value = 70272
$ws.Range("N85").Value = -72924
$ws.Range("H112").Value = 39129
$ws.Range("J112").Value = 39129
$ws.Range("L112").Value = 39129
$ws.Range("N112").Value = -42083
$ws.Range("H125").Value = 99999
$ws.Range("J125").Value = 99999
$ws.Range("L125").Value = 99999
$ws.Range("N125").Value = -109839
$ws.Range("H132").Value = 250051000
$ws.Range("I132").Value = 68002.336
$ws.Range("K132").Value = 204007.008
$ws.Range("M132").Value = -201477.008
$ws.Range("H134").Value = 47500
$ws.Range("J134").Value = 47500
$ws.Range("L134").Value = 142500
$ws.Range("N134").Value = -147570
$ws.Range("H136").Value = 5313.1333
$ws.Range("I136").Value = 5849.577
$ws.Range("K136").Value = 17548.731
$ws.Range("M136").Value = -14998.731
